# Auto-generated edit script applying the Siren_Profits.xlsx cell-value diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 417.23077
$ws.Range("I19").Value = 173.33333
$ws.Range("J19").Value = 490.4
$ws.Range("K19").Value = 173.33333
$ws.Range("L19").Value = 490.4
$ws.Range("M19").Value = 1.666670000000011
$ws.Range("N19").Value = -840.4
$ws.Range("H33").Value = 606.3125
$ws.Range("I33").Value = 353.92307
$ws.Range("K33").Value = 353.92307
$ws.Range("M33").Value = -124.92307
$ws.Range("H64").Value = 88927.86
$ws.Range("I64").Value = 150625
$ws.Range("J64").Value = 6665
$ws.Range("K64").Value = 150625
$ws.Range("L64").Value = 6665
$ws.Range("M64").Value = -150377
$ws.Range("N64").Value = -7161
$ws.Range("H67").Value = 88927.86
$ws.Range("I67").Value = 150625
$ws.Range("J67").Value = 6665
$ws.Range("K67").Value = 150625
$ws.Range("L67").Value = 6665
$ws.Range("M67").Value = -149767
$ws.Range("N67").Value = -8381
$ws.Range("H74").Value = 6312.5
$ws.Range("I74").Value = 6000
$ws.Range("J74").Value = 7250
$ws.Range("K74").Value = 6000
$ws.Range("L74").Value = 7250
$ws.Range("M74").Value = -5064
$ws.Range("N74").Value = -9122
$ws.Range("H77").Value = 6312.5
$ws.Range("I77").Value = 6000
$ws.Range("J77").Value = 7250
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 36250
$ws.Range("M77").Value = -25320
$ws.Range("N77").Value = -45610
$ws.Range("H80").Value = 50942.367
$ws.Range("I80").Value = 81082.28
$ws.Range("J80").Value = 3848.75
$ws.Range("K80").Value = 243246.84
$ws.Range("L80").Value = 11546.25
$ws.Range("M80").Value = -242248.84
$ws.Range("N80").Value = -13542.25
$ws.Range("H83").Value = 50942.367
$ws.Range("I83").Value = 81082.28
$ws.Range("J83").Value = 3848.75
$ws.Range("K83").Value = 729740.52
$ws.Range("L83").Value = 34638.75
$ws.Range("M83").Value = -724748.52
$ws.Range("N83").Value = -44622.75
$ws.Range("H88").Value = 1028.2858
$ws.Range("I88").Value = 649.75
$ws.Range("J88").Value = 1179.7
$ws.Range("K88").Value = 649.75
$ws.Range("L88").Value = 1179.7
$ws.Range("M88").Value = -243.75
$ws.Range("N88").Value = -1991.7
$ws.Range("H91").Value = 1028.2858
$ws.Range("I91").Value = 649.75
$ws.Range("J91").Value = 1179.7
$ws.Range("K91").Value = 649.75
$ws.Range("L91").Value = 1179.7
$ws.Range("M91").Value = 754.25
$ws.Range("N91").Value = -3987.7
$ws.Range("H98").Value = 21316.785
$ws.Range("I98").Value = 22802.73
$ws.Range("K98").Value = 22802.73
$ws.Range("M98").Value = -21304.73
$ws.Range("H103").Value = 798348.9
$ws.Range("I103").Value = 1529111.4
$ws.Range("J103").Value = 1153.4546
$ws.Range("K103").Value = 4587334.199999999
$ws.Range("L103").Value = 3460.3638
$ws.Range("M103").Value = -4586748.199999999
$ws.Range("N103").Value = -4632.3638
$ws.Range("H112").Value = 2149.1333
$ws.Range("J112").Value = 1936.125
$ws.Range("L112").Value = 5808.375
$ws.Range("N112").Value = -8024.375
$ws.Range("H116").Value = 5850733.5
$ws.Range("I116").Value = 9261296
$ws.Range("J116").Value = 4055.5715
$ws.Range("K116").Value = 9261296
$ws.Range("L116").Value = 4055.5715
$ws.Range("M116").Value = -9257854
$ws.Range("N116").Value = -10939.5715
$ws.Range("H122").Value = 21316.785
$ws.Range("I122").Value = 22802.73
$ws.Range("K122").Value = 68408.19
$ws.Range("M122").Value = -65958.19
$ws.Range("H132").Value = 1726948.8
$ws.Range("I132").Value = 2867.1345
$ws.Range("J132").Value = 16668989
$ws.Range("K132").Value = 8601.4035
$ws.Range("L132").Value = 50006967
$ws.Range("M132").Value = -6071.4035
$ws.Range("N132").Value = -50012027
$ws.Range("H135").Value = 7539.4165
$ws.Range("I135").Value = 10560.714
$ws.Range("J135").Value = 3309.6
$ws.Range("K135").Value = 95046.42600000001
$ws.Range("L135").Value = 29786.4
$ws.Range("M135").Value = -92511.42600000001
$ws.Range("N135").Value = -34856.39999999999
$ws.Range("H137").Value = 7895.636
$ws.Range("I137").Value = 10449.267
$ws.Range("J137").Value = 2423.5715
$ws.Range("K137").Value = 31347.801
$ws.Range("L137").Value = 7270.7145
$ws.Range("M137").Value = -28797.801
$ws.Range("N137").Value = -12370.7145
$ws.Range("H141").Value = 7304.4287
$ws.Range("I141").Value = 7357
$ws.Range("J141").Value = 6866.3335
$ws.Range("K141").Value = 22071
$ws.Range("L141").Value = 20599.0005
$ws.Range("M141").Value = -16891
$ws.Range("N141").Value = -30959.0005

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5850.518
$ws.Range("I32").Value = 5677.9634
$ws.Range("K32").Value = 5677.9634
$ws.Range("M32").Value = -5390.9634
$ws.Range("H51").Value = 60000
$ws.Range("J51").Value = 60000
$ws.Range("L51").Value = 60000
$ws.Range("N51").Value = -61512
$ws.Range("H61").Value = 8580.517
$ws.Range("I61").Value = 9699.9
$ws.Range("K61").Value = 9699.9
$ws.Range("M61").Value = -9487.9
$ws.Range("H74").Value = 1723.9615
$ws.Range("I74").Value = 941.4666999999999
$ws.Range("J74").Value = 2791
$ws.Range("K74").Value = 941.4666999999999
$ws.Range("L74").Value = 2791
$ws.Range("M74").Value = -67.46669999999995
$ws.Range("N74").Value = -4539
$ws.Range("H77").Value = 1723.9615
$ws.Range("I77").Value = 941.4666999999999
$ws.Range("J77").Value = 2791
$ws.Range("K77").Value = 4707.3335
$ws.Range("L77").Value = 13955
$ws.Range("M77").Value = -339.3334999999997
$ws.Range("N77").Value = -22691
$ws.Range("H132").Value = 3493.5
$ws.Range("I132").Value = 1555.2858
$ws.Range("J132").Value = 4726.909
$ws.Range("K132").Value = 4665.857400000001
$ws.Range("L132").Value = 14180.727
$ws.Range("M132").Value = -2135.857400000001
$ws.Range("N132").Value = -19240.727
$ws.Range("H135").Value = 136597.6
$ws.Range("J135").Value = 136597.6
$ws.Range("L135").Value = 136597.6
$ws.Range("N135").Value = -146737.6
$ws.Range("H136").Value = 8580.517
$ws.Range("I136").Value = 9699.9
$ws.Range("K136").Value = 29099.7
$ws.Range("M136").Value = -26549.7

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = ""
$ws.Range("H80").Value = 260.6
$ws.Range("J80").Value = 255.23077
$ws.Range("L80").Value = 255.23077
$ws.Range("N80").Value = -2251.23077
$ws.Range("H83").Value = 260.6
$ws.Range("J83").Value = 255.23077
$ws.Range("L83").Value = 1276.15385
$ws.Range("N83").Value = -11260.15385
$ws.Range("H86").Value = 5036.68
$ws.Range("I86").Value = 5548.45
$ws.Range("J86").Value = 2989.6
$ws.Range("K86").Value = 5548.45
$ws.Range("L86").Value = 2989.6
$ws.Range("M86").Value = -4425.45
$ws.Range("N86").Value = -5235.6
$ws.Range("H89").Value = 5036.68
$ws.Range("I89").Value = 5548.45
$ws.Range("J89").Value = 2989.6
$ws.Range("K89").Value = 27742.25
$ws.Range("L89").Value = 14948
$ws.Range("M89").Value = -22126.25
$ws.Range("N89").Value = -26180
$ws.Range("H107").Value = 2184.8667
$ws.Range("I107").Value = 2161.2727
$ws.Range("J107").Value = 2249.75
$ws.Range("K107").Value = 2161.2727
$ws.Range("L107").Value = 2249.75
$ws.Range("M107").Value = -241.2727
$ws.Range("N107").Value = -6089.75
$ws.Range("H134").Value = 5311.1665
$ws.Range("I134").Value = 5468.923
$ws.Range("K134").Value = 16406.769
$ws.Range("M134").Value = -13871.769

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 45499.5
$ws.Range("J54").Value = 45499.5
$ws.Range("L54").Value = 45499.5
$ws.Range("N54").Value = -46815.5
$ws.Range("H58").Value = 4745.75
$ws.Range("I58").Value = 5400.2
$ws.Range("K58").Value = 5400.2
$ws.Range("M58").Value = -5197.2
$ws.Range("H94").Value = 1045.8462
$ws.Range("I94").Value = 399.66666
$ws.Range("J94").Value = 1239.7
$ws.Range("K94").Value = 399.66666
$ws.Range("L94").Value = 1239.7
$ws.Range("M94").Value = 51.33334000000002
$ws.Range("N94").Value = -2141.7
$ws.Range("H99").Value = 266746.4
$ws.Range("I99").Value = 386860.16
$ws.Range("J99").Value = 6500
$ws.Range("K99").Value = 386860.16
$ws.Range("L99").Value = 6500
$ws.Range("M99").Value = -385362.16
$ws.Range("N99").Value = -9496
$ws.Range("H126").Value = 266746.4
$ws.Range("I126").Value = 386860.16
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 1160580.48
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -1158110.48
$ws.Range("N126").Value = -24440
$ws.Range("H132").Value = 1401.8387
$ws.Range("I132").Value = 1275.4615
$ws.Range("K132").Value = 3826.3845
$ws.Range("M132").Value = -1296.3845
$ws.Range("H134").Value = 4415.905
$ws.Range("I134").Value = 3412.7693
$ws.Range("J134").Value = 6046
$ws.Range("K134").Value = 10238.3079
$ws.Range("L134").Value = 18138
$ws.Range("M134").Value = -7703.3079
$ws.Range("N134").Value = -23208
$ws.Range("H135").Value = 86018
$ws.Range("J135").Value = 86018
$ws.Range("L135").Value = 86018
$ws.Range("N135").Value = -96158
$ws.Range("H136").Value = 4745.75
$ws.Range("I136").Value = 5400.2
$ws.Range("K136").Value = 16200.6
$ws.Range("M136").Value = -13650.6
$ws.Range("H141").Value = 194033.4
$ws.Range("J141").Value = 197748.89
$ws.Range("L141").Value = 197748.89
$ws.Range("N141").Value = -208108.89

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 357.8
$ws.Range("I8").Value = 357.8
$ws.Range("K8").Value = 1073.4
$ws.Range("M8").Value = -934.4000000000001
$ws.Range("H14").Value = 583.8333
$ws.Range("I14").Value = 583.8333
$ws.Range("K14").Value = 1751.4999
$ws.Range("M14").Value = -1578.4999
$ws.Range("H33").Value = 282.66666
$ws.Range("I33").Value = 274
$ws.Range("J33").Value = 285.14285
$ws.Range("K33").Value = 1644
$ws.Range("L33").Value = 1710.8571
$ws.Range("M33").Value = -1361
$ws.Range("N33").Value = -2276.8571
$ws.Range("H51").Value = 899.6667
$ws.Range("I51").Value = 599.5
$ws.Range("K51").Value = 1798.5
$ws.Range("M51").Value = -1338.5
$ws.Range("H131").Value = 10419304
$ws.Range("I131").Value = 71434970
$ws.Range("J131").Value = 1994.9878
$ws.Range("K131").Value = 214304910
$ws.Range("L131").Value = 5984.963400000001
$ws.Range("M131").Value = -214299870
$ws.Range("N131").Value = -16064.9634

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 18599.8
$ws.Range("J32").Value = 18599.8
$ws.Range("L32").Value = 18599.8
$ws.Range("N32").Value = -19191.8
$ws.Range("H97").Value = 6869.88
$ws.Range("I97").Value = 6715.8184
$ws.Range("J97").Value = 7999.6665
$ws.Range("K97").Value = 6715.8184
$ws.Range("L97").Value = 7999.6665
$ws.Range("M97").Value = -6219.8184
$ws.Range("N97").Value = -8991.666499999999
$ws.Range("H102").Value = 10713.286
$ws.Range("I102").Value = 12771.728
$ws.Range("K102").Value = 12771.728
$ws.Range("M102").Value = -11149.728
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""
$ws.Range("H122").Value = 10860.565
$ws.Range("I122").Value = 8427.556
$ws.Range("J122").Value = 19619.4
$ws.Range("K122").Value = 25282.668
$ws.Range("L122").Value = 58858.2
$ws.Range("M122").Value = -22832.668
$ws.Range("N122").Value = -63758.2
$ws.Range("H126").Value = 8290.477000000001
$ws.Range("I126").Value = 13319.4
$ws.Range("J126").Value = 3718.7273
$ws.Range("K126").Value = 39958.2
$ws.Range("L126").Value = 11156.1819
$ws.Range("M126").Value = -37488.2
$ws.Range("N126").Value = -16096.1819
$ws.Range("H132").Value = 1816.7872
$ws.Range("I132").Value = 1846.619
$ws.Range("J132").Value = 1566.2
$ws.Range("K132").Value = 5539.857
$ws.Range("L132").Value = 4698.6
$ws.Range("M132").Value = -3009.857
$ws.Range("N132").Value = -9758.6

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19278.5
$ws.Range("I7").Value = 22000.965
$ws.Range("K7").Value = 22000.965
$ws.Range("M7").Value = -21888.965
$ws.Range("H40").Value = 23410.408
$ws.Range("I40").Value = 31768
$ws.Range("K40").Value = 31768
$ws.Range("M40").Value = -31632
$ws.Range("H46").Value = 2011.4073
$ws.Range("I46").Value = 663.25
$ws.Range("J46").Value = 3089.9333
$ws.Range("K46").Value = 663.25
$ws.Range("L46").Value = 3089.9333
$ws.Range("M46").Value = -475.25
$ws.Range("N46").Value = -3465.9333
$ws.Range("H61").Value = 3840.88
$ws.Range("I61").Value = 1527.579
$ws.Range("K61").Value = 1527.579
$ws.Range("M61").Value = -1325.579
$ws.Range("H93").Value = 3581.2173
$ws.Range("I93").Value = 3935.4736
$ws.Range("K93").Value = 3935.4736
$ws.Range("M93").Value = -2687.4736
$ws.Range("H101").Value = 39674.25
$ws.Range("J101").Value = 39674.25
$ws.Range("L101").Value = 39674.25
$ws.Range("N101").Value = -46164.25
$ws.Range("H113").Value = 3840.88
$ws.Range("I113").Value = 1527.579
$ws.Range("K113").Value = 1527.579
$ws.Range("M113").Value = 642.421
$ws.Range("H122").Value = 7190.6113
$ws.Range("I122").Value = 6818
$ws.Range("J122").Value = 8159.4
$ws.Range("K122").Value = 20454
$ws.Range("L122").Value = 24478.2
$ws.Range("M122").Value = -18004
$ws.Range("N122").Value = -29378.2
$ws.Range("H126").Value = 19278.5
$ws.Range("I126").Value = 22000.965
$ws.Range("K126").Value = 66002.895
$ws.Range("M126").Value = -63532.895
$ws.Range("H132").Value = 499275.4
$ws.Range("I132").Value = 994619.6
$ws.Range("K132").Value = 2983858.8
$ws.Range("M132").Value = -2981328.8
$ws.Range("H141").Value = 89978.5
$ws.Range("J141").Value = 89978.5
$ws.Range("L141").Value = 89978.5
$ws.Range("N141").Value = -100338.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 99756.60000000001
$ws.Range("I74").Value = 19950
$ws.Range("K74").Value = 19950
$ws.Range("M74").Value = -19014
$ws.Range("H77").Value = 99756.60000000001
$ws.Range("I77").Value = 19950
$ws.Range("K77").Value = 59850
$ws.Range("M77").Value = -55170
$ws.Range("H107").Value = 24837.23
$ws.Range("I107").Value = 3061.8333
$ws.Range("J107").Value = 43501.855
$ws.Range("K107").Value = 9185.499899999999
$ws.Range("L107").Value = 130505.565
$ws.Range("M107").Value = -7265.499899999999
$ws.Range("N107").Value = -134345.565
$ws.Range("H113").Value = 1209.9286
$ws.Range("I113").Value = 743.1613
$ws.Range("K113").Value = 2229.4839
$ws.Range("M113").Value = -59.48390000000018
$ws.Range("H122").Value = 14996.047
$ws.Range("I122").Value = 2055.9312
$ws.Range("J122").Value = 41800.57
$ws.Range("K122").Value = 6167.7936
$ws.Range("L122").Value = 125401.71
$ws.Range("M122").Value = -3717.7936
$ws.Range("N122").Value = -130301.71
$ws.Range("H126").Value = 24140.953
$ws.Range("I126").Value = 36939.54
$ws.Range("J126").Value = 3343.25
$ws.Range("K126").Value = 110818.62
$ws.Range("L126").Value = 10029.75
$ws.Range("M126").Value = -108348.62
$ws.Range("N126").Value = -14969.75
$ws.Range("H132").Value = 9471.02
$ws.Range("I132").Value = 10899.919
$ws.Range("J132").Value = 5946.4
$ws.Range("K132").Value = 32699.757
$ws.Range("L132").Value = 17839.2
$ws.Range("M132").Value = -30169.757
$ws.Range("N132").Value = -22899.2
$ws.Range("H136").Value = 483187.3
$ws.Range("I136").Value = 594042
$ws.Range("J136").Value = 2817
$ws.Range("K136").Value = 1782126
$ws.Range("L136").Value = 8451
$ws.Range("M136").Value = -1779576
$ws.Range("N136").Value = -13551
